$wb = $excel.ActiveWorkbook

# Leave the previously-active sheet (Employee_Details) parked on A4,
# matching where the author last clicked before moving to the new sheet.
[void]$wb.Worksheets.Item("Employee_Details").Range("A4").Select()

# Remove the unwanted "Termination_Approval" feature file/sheet entirely...
[void]$wb.Worksheets.Item("Termination_Approval").Delete()

# ...and add a brand-new sheet in its place for the Line Manager
# "Change Salary" scenario (GL.AC.01).
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$ws.Name = "Line_Manager"

# Header row
$ws.Range("A1").Value = "scenario"
$ws.Range("B1").Value = "userName"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "employeeName"
$ws.Range("E1").Value = "actionManageSalary"
$ws.Range("F1").Value = "actionReasonManageSalary"
$ws.Range("G1").Value = "salaryAmount"

# Data row
$ws.Range("A2").Value = "LINEMANAGER_CHANGE_OF_SALARY_ACTION"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3056871"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "Welcome1"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "TestAutoFname TestAutoLname"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "Change Salary"

$ws.Range("F2").Value = "Change in Working Hours"

$ws.Range("G2").NumberFormat = "#,##0.00"
$ws.Range("G2").Value = 80572.46

# Make this newly-populated sheet the active tab/selection, like the author
# left it after editing.
[void]$ws.Activate()
[void]$ws.Range("A1:G2").Select()
